$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50 - this shifts existing rows 50-107 down to 51-108,
# carrying over the row's formatting (including the date-style D column).
$ws.Rows("50:50").Insert()

# Populate the newly inserted row 50 with the new weekly record.
$ws.Range("A50").Value = 8
$ws.Range("B50").Value = "Terminal La Palmera de La Serena"
$ws.Range("C50").Value = "Coquimbo"
$ws.Range("D50").Value = 44484
$ws.Range("E50").Value = 4
$ws.Range("F50").Value = 100112021
$ws.Range("G50").Value = "Ají"
$ws.Range("H50").Value = "Inferno"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 540
$ws.Range("K50").Value = 40000
$ws.Range("L50").Value = 41000
$ws.Range("M50").Value = 40500
$ws.Range("N50").Value = "$/caja 12 kilos"
$ws.Range("O50").Value = "Región de Arica y Parinacota"
$ws.Range("P50").Value = 3375
$ws.Range("Q50").Value = 12
$ws.Range("R50").Value = "Hortaliza"
